# Update the "Förändrad" (Changed) date column C for rows 2-33
# from serial date 46081 (2026-02-28) to 46082 (2026-03-01).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 46081) {
        $cell.Value = 46082
    }
}
